$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.838.34"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.876.84"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4603"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3881"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07872"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9858"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "1.947.80"
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.010"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.665"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06957"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009968"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "28.851.52"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.263"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.081"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.985"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.932"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09342"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9043"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.266"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.324"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.264"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.186"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05766"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.002"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.677"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1769"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.688"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.266"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5351"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07034"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.508"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.061"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
